$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.950.97'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.89%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.508.33'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -3.30%  '

$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '608.83'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.79'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -5.53%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.508.76'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.23%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.11%  '

$ws.Range("E9").Value = '  -0.95%  '

$ws.Range("E10").Value = '  -0.76%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.54'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.96%  '

$ws.Range("E12").Value = '  -1.91%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000216'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -4.30%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.01'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.50%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.098.81'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -3.49%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.519.21'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -3.12%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.991.79'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -1.83%  '

$ws.Range("E18").Value = '  -0.29%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.51'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.41%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.39'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.93%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.00'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '445.99'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.81%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.628'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.00'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.66%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.642.72'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -3.56%  '

$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("E27").Value = '  -8.83%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.99'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -6.00%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.69'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -4.96%  '

$ws.Range("E30").Value = '  -4.12%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.66'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -2.75%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.169'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.32%  '

$ws.Range("E33").Value = '  -0.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.63'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.85%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.19'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -6.20%  '

$ws.Range("E36").Value = '  -5.88%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.499.14'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.33%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.02'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.02%  '

$ws.Range("E39").Value = '  -0.01%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.28'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -3.05%  '

$ws.Range("E41").Value = '  -0.18%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '176.74'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.77%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0904'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.28%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.45'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.94%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '30.80'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -2.34%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.899'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.66%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.95'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.42%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.30'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -5.60%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.64'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.86%  '

$ws.Range("E50").Value = '  -11.21%  '

$ws.Range("E51").Value = '  -1.86%  '
